$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "24.869.34"
Set-TextValue $ws.Range("E2") "  -4.23%  "
Set-TextValue $ws.Range("D3") "1.631.67"
Set-TextValue $ws.Range("E3") "  -6.60%  "
Set-TextValue $ws.Range("D4") "0.9984"
Set-TextValue $ws.Range("E4") "  -0.13%  "
Set-TextValue $ws.Range("D5") "230.95"
Set-TextValue $ws.Range("E5") "  -7.42%  "
Set-TextValue $ws.Range("D6") "1.000"
Set-TextValue $ws.Range("E6") "  +0.07%  "
Set-TextValue $ws.Range("D7") "0.4724"
Set-TextValue $ws.Range("E7") "  -6.58%  "
Set-TextValue $ws.Range("D8") "0.2554"
Set-TextValue $ws.Range("E8") "  -7.17%  "
Set-TextValue $ws.Range("D9") "0.06062"
Set-TextValue $ws.Range("E9") "  -2.11%  "
Set-TextValue $ws.Range("D10") "0.06975"
Set-TextValue $ws.Range("E10") "  -3.95%  "
Set-TextValue $ws.Range("D11") "1.642.39"
Set-TextValue $ws.Range("E11") "  -5.94%  "
Set-TextValue $ws.Range("D12") "14.48"
Set-TextValue $ws.Range("E12") "  -4.74%  "
Set-TextValue $ws.Range("D13") "0.5930"
Set-TextValue $ws.Range("E13") "  -9.44%  "
Set-TextValue $ws.Range("D14") "4.320"
Set-TextValue $ws.Range("E14") "  -7.37%  "
Set-TextValue $ws.Range("D15") "73.06"
Set-TextValue $ws.Range("E15") "  -6.04%  "
Set-TextValue $ws.Range("D16") "0.9999"
Set-TextValue $ws.Range("E16") "  +0.03%  "
Set-TextValue $ws.Range("D17") "0.9995"
Set-TextValue $ws.Range("E17") "  +0.03%  "
Set-TextValue $ws.Range("D18") "24.871.09"
Set-TextValue $ws.Range("D19") "0.000006535"
Set-TextValue $ws.Range("E19") "  -4.69%  "
Set-TextValue $ws.Range("E20") "  -6.56%  "
Set-TextValue $ws.Range("D21") "1.851.11"
Set-TextValue $ws.Range("E21") "  -5.93%  "
Set-TextValue $ws.Range("D22") "4.312"
Set-TextValue $ws.Range("E22") "  -3.17%  "
Set-TextValue $ws.Range("D23") "8.486"
Set-TextValue $ws.Range("E23") "  -2.90%  "
Set-TextValue $ws.Range("D24") "5.184"
Set-TextValue $ws.Range("E24") "  -4.00%  "
Set-TextValue $ws.Range("D25") "132.58"
Set-TextValue $ws.Range("E25") "  -2.99%  "
Set-TextValue $ws.Range("D26") "14.71"
Set-TextValue $ws.Range("E26") "  -3.55%  "
Set-TextValue $ws.Range("D27") "1.381"
Set-TextValue $ws.Range("E27") "  -8.46%  "
Set-TextValue $ws.Range("D28") "102.85"
Set-TextValue $ws.Range("E28") "  -2.84%  "
Set-TextValue $ws.Range("D29") "1.621"
Set-TextValue $ws.Range("E29") "  -9.14%  "
Set-TextValue $ws.Range("D30") "3.865"
Set-TextValue $ws.Range("E30") "  -0.32%  "
Set-TextValue $ws.Range("D31") "0.07632"
Set-TextValue $ws.Range("E31") "  -6.95%  "
Set-TextValue $ws.Range("D32") "3.508"
Set-TextValue $ws.Range("E32") "  -3.88%  "
Set-TextValue $ws.Range("D33") "0.9993"
Set-TextValue $ws.Range("E33") "  +0.06%  "
Set-TextValue $ws.Range("D34") "0.04286"
Set-TextValue $ws.Range("E34") "  -8.52%  "
Set-TextValue $ws.Range("D35") "2.573"
Set-TextValue $ws.Range("E35") "  -3.05%  "
Set-TextValue $ws.Range("D36") "0.9168"
Set-TextValue $ws.Range("E36") "  -8.18%  "
Set-TextValue $ws.Range("D37") "0.5744"
Set-TextValue $ws.Range("E37") "  -6.98%  "
Set-TextValue $ws.Range("D38") "2.519"
Set-TextValue $ws.Range("E38") "  -8.34%  "
Set-TextValue $ws.Range("D41") "0.8080"
Set-TextValue $ws.Range("E41") "  +5.76%  "
Set-TextValue $ws.Range("D42") "96.89"
Set-TextValue $ws.Range("E42") "  -4.09%  "
Set-TextValue $ws.Range("D43") "1.740"
Set-TextValue $ws.Range("E43") "  -9.65%  "
Set-TextValue $ws.Range("D44") "0.3653"
Set-TextValue $ws.Range("E44") "  -6.98%  "
Set-TextValue $ws.Range("D45") "4.686"
Set-TextValue $ws.Range("E45") "  -6.41%  "
Set-TextValue $ws.Range("D46") "0.05181"
Set-TextValue $ws.Range("E46") "  -2.49%  "
Set-TextValue $ws.Range("D47") "0.1082"
Set-TextValue $ws.Range("E47") "  -6.06%  "
Set-TextValue $ws.Range("D48") "5.979"
Set-TextValue $ws.Range("E48") "  -5.72%  "
Set-TextValue $ws.Range("D49") "29.24"
Set-TextValue $ws.Range("E49") "  -4.82%  "
Set-TextValue $ws.Range("D50") "0.9995"
Set-TextValue $ws.Range("E50") "  -0.27%  "
Set-TextValue $ws.Range("D51") "0.9960"
Set-TextValue $ws.Range("E51") "  -0.37%  "

# Row 39/40 content swap: VeChain <-> PaxDollar
Set-TextValue $ws.Range("B39") "PaxDollar"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D39") "0.9982"
Set-TextValue $ws.Range("E39") "  -0.10%  "
Set-TextValue $ws.Range("B40") "VeChain"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D40") "0.01519"
Set-TextValue $ws.Range("E40") "  -6.04%  "
